$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column B (date strings like "2025-12-16") to be stored as text,
# matching the source data which keeps dates as literal strings rather than
# being auto-converted into Excel date serials.
$ws.Range("B2:B7").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'Thai League 1'
$ws.Range("B2").Value = '2025-12-16'
$ws.Range("C2").Value = '09:00:00'
$ws.Range("D2").Value = 'Rayong FC'
$ws.Range("E2").Value = 'Ratchaburi'
$ws.Range("F2").Value = 2.62
$ws.Range("G2").Value = 3.05
$ws.Range("H2").Value = 2.48
$ws.Range("I2").Value = 2.86
$ws.Range("J2").Value = 3.5
$ws.Range("K2").Value = 4.6
$ws.Range("L2").Value = 1.29
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 4.1
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.74
$ws.Range("R2").Value = 1.43
$ws.Range("S2").Value = 2.88
$ws.Range("T2").Value = 1.62
$ws.Range("U2").Value = 2.3
$ws.Range("V2").Value = 1.54
$ws.Range("W2").Value = 1.5
$ws.Range("X2").Value = 22
$ws.Range("Y2").Value = 15.5
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 46
$ws.Range("AB2").Value = 16.5
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 15
$ws.Range("AE2").Value = 34
$ws.Range("AF2").Value = 25
$ws.Range("AG2").Value = 15.5
$ws.Range("AH2").Value = 19.5
$ws.Range("AI2").Value = 44
$ws.Range("AJ2").Value = 55
$ws.Range("AK2").Value = 36
$ws.Range("AL2").Value = 46
$ws.Range("AM2").Value = 90
$ws.Range("AN2").Value = 27
$ws.Range("AO2").Value = 24

# Row 3
$ws.Range("A3").Value = 'Portuguese Segunda Liga'
$ws.Range("B3").Value = '2025-12-16'
$ws.Range("C3").Value = '14:00:00'
$ws.Range("D3").Value = 'Maritimo'
$ws.Range("E3").Value = 'Benfica B'
$ws.Range("F3").Value = 1.91
$ws.Range("G3").Value = 1.98
$ws.Range("H3").Value = 4.2
$ws.Range("I3").Value = 4.7
$ws.Range("J3").Value = 3.6
$ws.Range("K3").Value = 4.1
$ws.Range("L3").Value = 1.41
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.2
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 1.92
$ws.Range("Q3").Value = 1.84
$ws.Range("R3").Value = 1.34
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1.78
$ws.Range("U3").Value = 2.04
$ws.Range("V3").Value = 1.27
$ws.Range("W3").Value = 2.02
$ws.Range("X3").Value = 16.5
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 10.5
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("A4").Value = 'Swiss Super League'
$ws.Range("B4").Value = '2025-12-16'
$ws.Range("C4").Value = '16:30:00'
$ws.Range("D4").Value = 'St Gallen'
$ws.Range("E4").Value = 'Sion'
$ws.Range("F4").Value = 2.32
$ws.Range("G4").Value = 2.58
$ws.Range("H4").Value = 2.92
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 3.65
$ws.Range("K4").Value = 4.1
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 2.2
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 1.68
$ws.Range("R4").Value = 1.4
$ws.Range("S4").Value = 2.42
$ws.Range("T4").Value = 1.5
$ws.Range("U4").Value = 1.01
$ws.Range("V4").Value = 1.41
$ws.Range("W4").Value = 1.63
$ws.Range("X4").Value = 26
$ws.Range("Y4").Value = 21
$ws.Range("Z4").Value = 34
$ws.Range("AA4").Value = 70
$ws.Range("AB4").Value = 18.5
$ws.Range("AC4").Value = 10.5
$ws.Range("AD4").Value = 20
$ws.Range("AE4").Value = 46
$ws.Range("AF4").Value = 25
$ws.Range("AG4").Value = 17.5
$ws.Range("AH4").Value = 22
$ws.Range("AI4").Value = 55
$ws.Range("AJ4").Value = 44
$ws.Range("AK4").Value = 32
$ws.Range("AL4").Value = 44
$ws.Range("AM4").Value = 100
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("A5").Value = 'Swiss Super League'
$ws.Range("B5").Value = '2025-12-16'
$ws.Range("C5").Value = '16:30:00'
$ws.Range("D5").Value = 'Winterthur'
$ws.Range("E5").Value = 'Thun'
$ws.Range("F5").Value = 4.1
$ws.Range("G5").Value = 4.6
$ws.Range("H5").Value = 1.77
$ws.Range("I5").Value = 1.87
$ws.Range("J5").Value = 4.3
$ws.Range("K5").Value = 4.9
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 6.2
$ws.Range("O5").Value = 1.16
$ws.Range("P5").Value = 2.8
$ws.Range("Q5").Value = 1.49
$ws.Range("R5").Value = 1.73
$ws.Range("S5").Value = 2.16
$ws.Range("T5").Value = 1.44
$ws.Range("U5").Value = 2.32
$ws.Range("V5").Value = 2.14
$ws.Range("W5").Value = 1.28
$ws.Range("X5").Value = 34
$ws.Range("Y5").Value = 15
$ws.Range("Z5").Value = 15.5
$ws.Range("AA5").Value = 21
$ws.Range("AB5").Value = 25
$ws.Range("AC5").Value = 11.5
$ws.Range("AD5").Value = 11.5
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 38
$ws.Range("AG5").Value = 18.5
$ws.Range("AH5").Value = 16
$ws.Range("AI5").Value = 25
$ws.Range("AJ5").Value = 80
$ws.Range("AK5").Value = 42
$ws.Range("AL5").Value = 40
$ws.Range("AM5").Value = 55
$ws.Range("AN5").Value = 29
$ws.Range("AO5").Value = 7

# Row 6
$ws.Range("A6").Value = 'English National League'
$ws.Range("B6").Value = '2025-12-16'
$ws.Range("C6").Value = '16:45:00'
$ws.Range("D6").Value = 'Truro City'
$ws.Range("E6").Value = 'Wealdstone'
$ws.Range("F6").Value = 2.7
$ws.Range("G6").Value = 3.05
$ws.Range("H6").Value = 2.48
$ws.Range("I6").Value = 2.8
$ws.Range("J6").Value = 3.35
$ws.Range("K6").Value = 3.95
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 3.9
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 1.98
$ws.Range("Q6").Value = 1.82
$ws.Range("R6").Value = 1.13
$ws.Range("S6").Value = 1.85
$ws.Range("T6").Value = 1.01
$ws.Range("U6").Value = 2.12
$ws.Range("V6").Value = 1.55
$ws.Range("W6").Value = 1.5
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 16.5
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 38
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 60
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 100
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000

# Row 7
$ws.Range("A7").Value = 'Welsh Premiership'
$ws.Range("B7").Value = '2025-12-16'
$ws.Range("C7").Value = '16:45:00'
$ws.Range("D7").Value = 'Cardiff Metropolitan'
$ws.Range("E7").Value = 'Briton Ferry Llansawel'
$ws.Range("F7").Value = 1.54
$ws.Range("G7").Value = 1.93
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 7.6
$ws.Range("J7").Value = 2.84
$ws.Range("K7").Value = 5
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 1.9
$ws.Range("O7").Value = 1.01
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.76
$ws.Range("R7").Value = 1.24
$ws.Range("S7").Value = 2.6
$ws.Range("T7").Value = 1.01
$ws.Range("U7").Value = 1.65
$ws.Range("V7").Value = 1.15
$ws.Range("W7").Value = 2.06
$ws.Range("X7").Value = 22
$ws.Range("Y7").Value = 28
$ws.Range("Z7").Value = 70
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 12
$ws.Range("AC7").Value = 13
$ws.Range("AD7").Value = 34
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 15
$ws.Range("AG7").Value = 14.5
$ws.Range("AH7").Value = 30
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 25
$ws.Range("AK7").Value = 28
$ws.Range("AL7").Value = 55
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000
